# Edit the "存款" (deposit) sheet: add header labels for columns B-E,
# add new columns F-M (total, property_category, category, date,
# legislator_name, legislator_id, source_file, index) and populate the
# `total` amount (column F) for every data row, shifting the previous
# (mis-placed) amount that lived in column G into its correct spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 2).Value2 = "bank"
$ws.Cells.Item(1, 3).Value2 = "deposit_type"
$ws.Cells.Item(1, 4).Value2 = "currency"
$ws.Cells.Item(1, 5).Value2 = "owner"
$ws.Cells.Item(1, 6).Value2 = "total"
$ws.Cells.Item(1, 7).Value2 = "property_category"
$ws.Cells.Item(1, 8).Value2 = "category"
$ws.Cells.Item(1, 9).Value2 = "date"
$ws.Cells.Item(1, 10).Value2 = "legislator_name"
$ws.Cells.Item(1, 11).Value2 = "legislator_id"
$ws.Cells.Item(1, 12).Value2 = "source_file"
$ws.Cells.Item(1, 13).Value2 = "index"

# ---- Data rows (2-17): the `total` amount values, keyed by row ----
$totals = @{
    2  = 1219105
    3  = 1176037
    4  = 1589557
    5  = 2937
    6  = 5380
    7  = 3059
    8  = 520466
    9  = 223262
    10 = 1278554
    11 = 200000
    12 = 491331
    13 = 29068
    14 = 456043
    15 = 568433
    16 = 462276.89
    17 = 34.52
}

for ($row = 2; $row -le 17; $row++) {
    $index = $ws.Cells.Item($row, 1).Value2

    # Column F: total amount (was missing / mis-placed in column G before)
    $ws.Cells.Item($row, 6).Value2 = $totals[$row]

    # Columns G-M: constant metadata columns shared by every deposit row
    $ws.Cells.Item($row, 7).Value2  = "deposit"
    $ws.Cells.Item($row, 8).Value2  = "normal"
    $ws.Cells.Item($row, 9).Value2  = "2012-04-30"
    $ws.Cells.Item($row, 10).Value2 = "楊麗環"
    $ws.Cells.Item($row, 11).Value2 = 960
    $ws.Cells.Item($row, 12).Value2 = "tmp700a1"
    $ws.Cells.Item($row, 13).Value2 = $index
}
